$wb = $excel.ActiveWorkbook

# --- Update selection/tab state on the existing "dockermetrics" sheet ---
$dockermetrics = $wb.Worksheets.Item("dockermetrics")
$dockermetrics.Range("A1:D1").Select()

# --- Add the new "packetbeats" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "packetbeats"

# --- Header row ---
$ws.Range("A1").Value = "Field"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Example"
$ws.Range("D1").Value = "Data Type"
$ws.Range("A1:D1").Font.Bold = $true

# --- Data rows (Field column) ---
$ws.Cells.Item(2, 1).Value = "@timestamp"
$ws.Cells.Item(3, 1).Value = "_id"
$ws.Cells.Item(4, 1).Value = "_index"
$ws.Cells.Item(5, 1).Value = "_score"
$ws.Cells.Item(6, 1).Value = "_type"
$ws.Cells.Item(7, 1).Value = "agent.ephemeral_id"
$ws.Cells.Item(8, 1).Value = "agent.hostname"
$ws.Cells.Item(9, 1).Value = "agent.id"
$ws.Cells.Item(10, 1).Value = "agent.type"
$ws.Cells.Item(11, 1).Value = "agent.version"
$ws.Cells.Item(12, 1).Value = "client.bytes"
$ws.Cells.Item(13, 1).Value = "client.ip"
$ws.Cells.Item(14, 1).Value = "client.port"
$ws.Cells.Item(15, 1).Value = "destination.bytes"
$ws.Cells.Item(16, 1).Value = "destination.domain"
$ws.Cells.Item(17, 1).Value = "destination.ip"
$ws.Cells.Item(18, 1).Value = "destination.packets"
$ws.Cells.Item(19, 1).Value = "destination.port"
$ws.Cells.Item(20, 1).Value = "dns.additionals_count"
$ws.Cells.Item(21, 1).Value = "dns.answers_count"
$ws.Cells.Item(22, 1).Value = "dns.authorities_count"
$ws.Cells.Item(23, 1).Value = "dns.flags.authentic_data"
$ws.Cells.Item(24, 1).Value = "dns.flags.authoritative"
$ws.Cells.Item(25, 1).Value = "dns.flags.checking_disabled"
$ws.Cells.Item(26, 1).Value = "dns.flags.recursion_available"
$ws.Cells.Item(27, 1).Value = "dns.flags.recursion_desired"
$ws.Cells.Item(28, 1).Value = "dns.flags.truncated_response"
$ws.Cells.Item(29, 1).Value = "dns.header_flags"
$ws.Cells.Item(30, 1).Value = "dns.id"
$ws.Cells.Item(31, 1).Value = "dns.op_code"
$ws.Cells.Item(32, 1).Value = "dns.question.class"
$ws.Cells.Item(33, 1).Value = "dns.question.etld_plus_one"
$ws.Cells.Item(34, 1).Value = "dns.question.name"
$ws.Cells.Item(35, 1).Value = "dns.question.registered_domain"
$ws.Cells.Item(36, 1).Value = "dns.question.type"
$ws.Cells.Item(37, 1).Value = "dns.resolved_ip"
$ws.Cells.Item(38, 1).Value = "dns.response_code"
$ws.Cells.Item(39, 1).Value = "dns.type"
$ws.Cells.Item(40, 1).Value = "ecs.version"
$ws.Cells.Item(41, 1).Value = "error.message"
$ws.Cells.Item(42, 1).Value = "event.action"
$ws.Cells.Item(43, 1).Value = "event.category"
$ws.Cells.Item(44, 1).Value = "event.dataset"
$ws.Cells.Item(45, 1).Value = "event.duration"
$ws.Cells.Item(46, 1).Value = "event.end"
$ws.Cells.Item(47, 1).Value = "event.kind"
$ws.Cells.Item(48, 1).Value = "event.start"
$ws.Cells.Item(49, 1).Value = "flow.final"
$ws.Cells.Item(50, 1).Value = "flow.id"
$ws.Cells.Item(51, 1).Value = "host.name"
$ws.Cells.Item(52, 1).Value = "http.request.body.bytes"
$ws.Cells.Item(53, 1).Value = "http.request.bytes"
$ws.Cells.Item(54, 1).Value = "http.request.headers.accept"
$ws.Cells.Item(55, 1).Value = "http.request.headers.accept-encoding"
$ws.Cells.Item(56, 1).Value = "http.request.headers.accept-language"
$ws.Cells.Item(57, 1).Value = "http.request.headers.authorization"
$ws.Cells.Item(58, 1).Value = "http.request.headers.connection"
$ws.Cells.Item(59, 1).Value = "http.request.headers.content-length"
$ws.Cells.Item(60, 1).Value = "http.request.headers.content-type"
$ws.Cells.Item(61, 1).Value = "http.request.headers.host"
$ws.Cells.Item(62, 1).Value = "http.request.headers.if-none-match"
$ws.Cells.Item(63, 1).Value = "http.request.headers.kbn-version"
$ws.Cells.Item(64, 1).Value = "http.request.headers.origin"
$ws.Cells.Item(65, 1).Value = "http.request.headers.referer"
$ws.Cells.Item(66, 1).Value = "http.request.headers.user-agent"
$ws.Cells.Item(67, 1).Value = "http.request.method"
$ws.Cells.Item(68, 1).Value = "http.request.referrer"
$ws.Cells.Item(69, 1).Value = "http.response.body.bytes"
$ws.Cells.Item(70, 1).Value = "http.response.bytes"
$ws.Cells.Item(71, 1).Value = "http.response.headers.accept-ranges"
$ws.Cells.Item(72, 1).Value = "http.response.headers.cache-control"
$ws.Cells.Item(73, 1).Value = "http.response.headers.connection"
$ws.Cells.Item(74, 1).Value = "http.response.headers.content-encoding"
$ws.Cells.Item(75, 1).Value = "http.response.headers.content-length"
$ws.Cells.Item(76, 1).Value = "http.response.headers.content-security-policy"
$ws.Cells.Item(77, 1).Value = "http.response.headers.content-type"
$ws.Cells.Item(78, 1).Value = "http.response.headers.date"
$ws.Cells.Item(79, 1).Value = "http.response.headers.etag"
$ws.Cells.Item(80, 1).Value = "http.response.headers.kbn-name"
$ws.Cells.Item(81, 1).Value = "http.response.headers.kbn-xpack-sig"
$ws.Cells.Item(82, 1).Value = "http.response.headers.location"
$ws.Cells.Item(83, 1).Value = "http.response.headers.transfer-encoding"
$ws.Cells.Item(84, 1).Value = "http.response.headers.vary"
$ws.Cells.Item(85, 1).Value = "http.response.status_code"
$ws.Cells.Item(86, 1).Value = "http.response.status_phrase"
$ws.Cells.Item(87, 1).Value = "http.version"
$ws.Cells.Item(88, 1).Value = "icmp.request.code"
$ws.Cells.Item(89, 1).Value = "icmp.request.message"
$ws.Cells.Item(90, 1).Value = "icmp.request.type"
$ws.Cells.Item(91, 1).Value = "icmp.response.code"
$ws.Cells.Item(92, 1).Value = "icmp.response.message"
$ws.Cells.Item(93, 1).Value = "icmp.response.type"
$ws.Cells.Item(94, 1).Value = "icmp.version"
$ws.Cells.Item(95, 1).Value = "method"
$ws.Cells.Item(96, 1).Value = "network.bytes"
$ws.Cells.Item(97, 1).Value = "network.community_id"
$ws.Cells.Item(98, 1).Value = "network.direction"
$ws.Cells.Item(99, 1).Value = "network.packets"
$ws.Cells.Item(100, 1).Value = "network.protocol"
$ws.Cells.Item(101, 1).Value = "network.transport"
$ws.Cells.Item(102, 1).Value = "network.type"
$ws.Cells.Item(103, 1).Value = "path"
$ws.Cells.Item(104, 1).Value = "query"
$ws.Cells.Item(105, 1).Value = "resource"
$ws.Cells.Item(106, 1).Value = "server.bytes"
$ws.Cells.Item(107, 1).Value = "server.domain"
$ws.Cells.Item(108, 1).Value = "server.ip"
$ws.Cells.Item(109, 1).Value = "server.port"
$ws.Cells.Item(110, 1).Value = "source.bytes"
$ws.Cells.Item(111, 1).Value = "source.ip"
$ws.Cells.Item(112, 1).Value = "source.packets"
$ws.Cells.Item(113, 1).Value = "source.port"
$ws.Cells.Item(114, 1).Value = "status"
$ws.Cells.Item(115, 1).Value = "type"
$ws.Cells.Item(116, 1).Value = "url.domain"
$ws.Cells.Item(117, 1).Value = "url.full"
$ws.Cells.Item(118, 1).Value = "url.path"
$ws.Cells.Item(119, 1).Value = "url.port"
$ws.Cells.Item(120, 1).Value = "url.query"
$ws.Cells.Item(121, 1).Value = "url.scheme"
$ws.Cells.Item(122, 1).Value = "user_agent.original"

$ws.Range("A2:A122").Font.Color = 0

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 38.333333333333336
$ws.Columns.Item(2).ColumnWidth = 21.666666666666668
$ws.Columns.Item(3).ColumnWidth = 35.166666666666664
$ws.Columns.Item(4).ColumnWidth = 25

# --- Selection on new sheet ---
$ws.Range("A6").Select()
